$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Pausas Activas" records appended below the existing data (rows 65-73).
# Column A holds numeric-looking employee IDs that must stay TEXT (matching the
# rest of the sheet, which stores every value as a string), so a text number
# format is forced before writing the value, then the default "Normal" style is
# restored so no stray formatting is left behind on the cell.
# Each record is wrapped with the unary "," so PowerShell keeps it as a nested
# array element instead of flattening it into the outer $data array.

$data = @(
    ,@(65, '11639', 'EDRIAN ALONSO BENITEZ ESCOBAR', '2025-05-16T10:52', 'SELL E', 'Sí', '10:55:12 a.m.', '10:55:59 a.m.', 46)
    ,@(66, '11639', 'EDRIAN ALONSO BENITEZ ESCOBAR', '2025-05-16T10:56', 'SELL E', 'No', '10:57:08 a.m.', '10:57:09 a.m.', 1)
    ,@(67, '11639', 'EDRIAN ALONSO BENITEZ ESCOBAR', '2025-05-20T14:32', 'SELL E', 'Sí', '2:33:39 p.m.', '2:34:38 p.m.', 59)
    ,@(68, '11227', 'JOHANA ANDREA GONZALEZ LOPEZ', '2025-05-20T14:34', 'ADIT E', 'Sí', '2:34:54 p.m.', '2:35:18 p.m.', 24)
    ,@(69, '11636', 'YOVANI ANDRES  BEDOYA PEREZ', '2025-05-26T14:38', 'SELL E', 'Sí', '2:38:50 p.m.', '2:39:32 p.m.', 41)
    ,@(70, '10640', 'Jhonnattan Ruiz', '2025-06-12T13:09', 'ADM E', 'Sí', '1:49:02 p.m.', '1:49:12 p.m.', 9)
    ,@(71, '11639', 'EDRIAN ALONSO BENITEZ ESCOBAR', '2025-07-09T13:52', 'SELL E', 'Sí', '1:52:53 p.m.', '1:53:05 p.m.', 12)
    ,@(72, '11639', 'EDRIAN ALONSO BENITEZ ESCOBAR', '2025-08-05T12:55', 'SELL E', 'Sí', '12:56:01 p.m.', '12:56:40 p.m.', 38)
    ,@(73, '11639', 'EDRIAN ALONSO BENITEZ ESCOBAR', '2025-08-05T13:20', 'SELL E', 'Sí', '1:20:13 p.m.', '1:22:13 p.m.', 120)
)

foreach ($row in $data) {
    $r = $row[0]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[1]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
